$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 106.2
$ws.Range("I2").Value = 106.2
$ws.Range("K2").Value = 106.2
$ws.Range("M2").Value = 6.799999999999997
$ws.Range("H6").Value = 81520.62
$ws.Range("I6").Value = 81520.62
$ws.Range("K6").Value = 244561.86
$ws.Range("M6").Value = -244449.86
$ws.Range("H11").Value = 97725.25
$ws.Range("I11").Value = 97725.25
$ws.Range("K11").Value = 97725.25
$ws.Range("M11").Value = -97585.25
$ws.Range("H12").Value = 167
$ws.Range("I12").Value = 167
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 167
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 3
$ws.Range("N12").ClearContents()
$ws.Range("H17").Value = 672.75757
$ws.Range("J17").Value = 611.1404
$ws.Range("L17").Value = 1833.4212
$ws.Range("N17").Value = -2169.4212
$ws.Range("H19").Value = 1210.5
$ws.Range("I19").Value = 1230
$ws.Range("J19").Value = 1198.8
$ws.Range("K19").Value = 1230
$ws.Range("L19").Value = 1198.8
$ws.Range("M19").Value = -1055
$ws.Range("N19").Value = -1548.8
$ws.Range("H28").Value = 482.07407
$ws.Range("I28").Value = 462.22726
$ws.Range("J28").Value = 569.4
$ws.Range("K28").Value = 462.22726
$ws.Range("L28").Value = 569.4
$ws.Range("M28").Value = 22.77274
$ws.Range("N28").Value = -1539.4
$ws.Range("H33").Value = 252.71428
$ws.Range("I33").Value = 233.6923
$ws.Range("K33").Value = 233.6923
$ws.Range("M33").Value = -4.692299999999989
$ws.Range("H40").Value = 5480
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H53").Value = 496.1111
$ws.Range("J53").Value = 766.25
$ws.Range("L53").Value = 766.25
$ws.Range("N53").Value = -2040.25
$ws.Range("H82").Value = 1104
$ws.Range("I82").Value = 1104
$ws.Range("K82").Value = 3312
$ws.Range("M82").Value = -2906
$ws.Range("H85").Value = 1104
$ws.Range("I85").Value = 1104
$ws.Range("K85").Value = 3312
$ws.Range("M85").Value = -1908
$ws.Range("H86").Value = 13525.125
$ws.Range("I86").Value = 14700.5
$ws.Range("J86").Value = 12349.75
$ws.Range("K86").Value = 14700.5
$ws.Range("L86").Value = 12349.75
$ws.Range("M86").Value = -13577.5
$ws.Range("N86").Value = -14595.75
$ws.Range("H89").Value = 13525.125
$ws.Range("I89").Value = 14700.5
$ws.Range("J89").Value = 12349.75
$ws.Range("K89").Value = 73502.5
$ws.Range("L89").Value = 61748.75
$ws.Range("M89").Value = -67886.5
$ws.Range("N89").Value = -72980.75
$ws.Range("H96").Value = 1921.8125
$ws.Range("I96").Value = 459.83334
$ws.Range("K96").Value = 1379.50002
$ws.Range("M96").Value = -6.50001999999995
$ws.Range("H98").Value = 710.8214
$ws.Range("I98").Value = 820.41174
$ws.Range("J98").Value = 541.4545000000001
$ws.Range("K98").Value = 820.41174
$ws.Range("L98").Value = 541.4545000000001
$ws.Range("M98").Value = 677.58826
$ws.Range("N98").Value = -3537.4545
$ws.Range("H99").Value = 188
$ws.Range("I99").Value = 188
$ws.Range("K99").Value = 564
$ws.Range("M99").Value = 934
$ws.Range("H112").Value = 2757
$ws.Range("J112").Value = 2656.4285
$ws.Range("L112").Value = 7969.2855
$ws.Range("N112").Value = -10185.2855
$ws.Range("H113").Value = 4207.9165
$ws.Range("I113").Value = 4221.8887
$ws.Range("J113").Value = 4166
$ws.Range("K113").Value = 4221.8887
$ws.Range("L113").Value = 4166
$ws.Range("M113").Value = -967.8887000000004
$ws.Range("N113").Value = -10674
$ws.Range("H118").Value = 1789.125
$ws.Range("J118").Value = 2624.75
$ws.Range("L118").Value = 7874.25
$ws.Range("N118").Value = -11188.25
$ws.Range("H121").Value = 4990
$ws.Range("J121").Value = 4990
$ws.Range("L121").Value = 14970
$ws.Range("N121").Value = -18464
$ws.Range("H122").Value = 710.8214
$ws.Range("I122").Value = 820.41174
$ws.Range("J122").Value = 541.4545000000001
$ws.Range("K122").Value = 2461.23522
$ws.Range("L122").Value = 1624.3635
$ws.Range("M122").Value = -11.23522000000003
$ws.Range("N122").Value = -6524.3635
$ws.Range("H132").Value = 17758.084
$ws.Range("I132").Value = 19300.092
$ws.Range("J132").Value = 796
$ws.Range("K132").Value = 57900.276
$ws.Range("L132").Value = 2388
$ws.Range("M132").Value = -55370.276
$ws.Range("N132").Value = -7448
$ws.Range("H135").Value = 5270.241
$ws.Range("I135").Value = 987.53845
$ws.Range("J135").Value = 8749.9375
$ws.Range("K135").Value = 8887.84605
$ws.Range("L135").Value = 78749.4375
$ws.Range("M135").Value = -6352.84605
$ws.Range("N135").Value = -83819.4375
$ws.Range("H137").Value = 3615.3877
$ws.Range("I137").Value = 1577.3864
$ws.Range("K137").Value = 4732.1592
$ws.Range("M137").Value = -2182.1592
$ws.Range("H138").Value = 3593.32
$ws.Range("I138").Value = 2835.1667
$ws.Range("K138").Value = 8505.500100000001
$ws.Range("M138").Value = -3365.500100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 736.6
$ws.Range("I2").Value = 736.4761999999999
$ws.Range("J2").Value = 737.25
$ws.Range("K2").Value = 736.4761999999999
$ws.Range("L2").Value = 737.25
$ws.Range("M2").Value = -623.4761999999999
$ws.Range("N2").Value = -963.25
$ws.Range("H32").Value = 223417.23
$ws.Range("I32").Value = 231208.94
$ws.Range("K32").Value = 231208.94
$ws.Range("M32").Value = -230921.94
$ws.Range("H45").Value = 2485.9167
$ws.Range("I45").Value = 1955.1666
$ws.Range("K45").Value = 1955.1666
$ws.Range("M45").Value = -1578.1666
$ws.Range("H61").Value = 3349.2307
$ws.Range("I61").Value = 3421.9092
$ws.Range("J61").Value = 2949.5
$ws.Range("K61").Value = 3421.9092
$ws.Range("L61").Value = 2949.5
$ws.Range("M61").Value = -3209.9092
$ws.Range("N61").Value = -3373.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H97").Value = 1859.75
$ws.Range("I97").Value = 1413
$ws.Range("K97").Value = 1413
$ws.Range("M97").Value = -917
$ws.Range("H102").Value = 1633.1333
$ws.Range("I102").Value = 1269.8462
$ws.Range("K102").Value = 1269.8462
$ws.Range("M102").Value = 352.1538
$ws.Range("H116").Value = 736.6
$ws.Range("I116").Value = 736.4761999999999
$ws.Range("J116").Value = 737.25
$ws.Range("K116").Value = 736.4761999999999
$ws.Range("L116").Value = 737.25
$ws.Range("M116").Value = 1557.5238
$ws.Range("N116").Value = -5325.25
$ws.Range("H122").Value = 2344.1396
$ws.Range("I122").Value = 1903.0938
$ws.Range("J122").Value = 3627.182
$ws.Range("K122").Value = 5709.2814
$ws.Range("L122").Value = 10881.546
$ws.Range("M122").Value = -3259.2814
$ws.Range("N122").Value = -15781.546
$ws.Range("H126").Value = 5249.5
$ws.Range("I126").Value = 5249.5
$ws.Range("K126").Value = 15748.5
$ws.Range("M126").Value = -13278.5
$ws.Range("H132").Value = 4579.5776
$ws.Range("I132").Value = 2752.8438
$ws.Range("K132").Value = 8258.5314
$ws.Range("M132").Value = -5728.5314
$ws.Range("H136").Value = 3349.2307
$ws.Range("I136").Value = 3421.9092
$ws.Range("J136").Value = 2949.5
$ws.Range("K136").Value = 10265.7276
$ws.Range("L136").Value = 8848.5
$ws.Range("M136").Value = -7715.7276
$ws.Range("N136").Value = -13948.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 736.6
$ws.Range("I3").Value = 736.4761999999999
$ws.Range("J3").Value = 737.25
$ws.Range("K3").Value = 736.4761999999999
$ws.Range("L3").Value = 737.25
$ws.Range("M3").Value = -622.4761999999999
$ws.Range("N3").Value = -965.25
$ws.Range("H22").Value = 261.33334
$ws.Range("I22").Value = 276.1
$ws.Range("J22").Value = 187.5
$ws.Range("K22").Value = 276.1
$ws.Range("L22").Value = 187.5
$ws.Range("M22").Value = -103.1
$ws.Range("N22").Value = -533.5
$ws.Range("H60").Value = 105000
$ws.Range("J60").Value = 105000
$ws.Range("L60").Value = 105000
$ws.Range("N60").Value = -106198
$ws.Range("H68").Value = 25000
$ws.Range("I68").Value = 25000
$ws.Range("K68").Value = 25000
$ws.Range("M68").Value = -24189
$ws.Range("H71").Value = 25000
$ws.Range("I71").Value = 25000
$ws.Range("K71").Value = 75000
$ws.Range("M71").Value = -70944
$ws.Range("H86").Value = 4380.1177
$ws.Range("I86").Value = 4461.5713
$ws.Range("K86").Value = 4461.5713
$ws.Range("M86").Value = -3338.5713
$ws.Range("H89").Value = 4380.1177
$ws.Range("I89").Value = 4461.5713
$ws.Range("K89").Value = 22307.8565
$ws.Range("M89").Value = -16691.8565
$ws.Range("H94").Value = 4719.9473
$ws.Range("I94").Value = 4542.4375
$ws.Range("J94").Value = 5666.6665
$ws.Range("K94").Value = 4542.4375
$ws.Range("L94").Value = 5666.6665
$ws.Range("M94").Value = -4091.4375
$ws.Range("N94").Value = -6568.6665
$ws.Range("H99").Value = 26500
$ws.Range("I99").Value = 26500
$ws.Range("K99").Value = 26500
$ws.Range("M99").Value = -25002
$ws.Range("H105").Value = 4944.5
$ws.Range("I105").Value = 1939.6
$ws.Range("K105").Value = 1939.6
$ws.Range("M105").Value = -192.5999999999999
$ws.Range("H107").Value = 1432.9445
$ws.Range("I107").Value = 1222.3529
$ws.Range("J107").Value = 5013
$ws.Range("K107").Value = 1222.3529
$ws.Range("L107").Value = 5013
$ws.Range("M107").Value = 697.6470999999999
$ws.Range("N107").Value = -8853
$ws.Range("H134").Value = 7075.7896
$ws.Range("I134").Value = 7465.4375
$ws.Range("J134").Value = 4997.6665
$ws.Range("K134").Value = 22396.3125
$ws.Range("L134").Value = 14992.9995
$ws.Range("M134").Value = -19861.3125
$ws.Range("N134").Value = -20062.9995
$ws.Range("H135").Value = 66491.2
$ws.Range("J135").Value = 66491.2
$ws.Range("L135").Value = 66491.2
$ws.Range("N135").Value = -76631.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 33424.668
$ws.Range("J7").Value = 132.83333
$ws.Range("L7").Value = 132.83333
$ws.Range("N7").Value = -358.83333
$ws.Range("H22").Value = 1104.84
$ws.Range("I22").Value = 647.4706
$ws.Range("J22").Value = 2076.75
$ws.Range("K22").Value = 647.4706
$ws.Range("L22").Value = 2076.75
$ws.Range("M22").Value = -297.4706
$ws.Range("N22").Value = -2776.75
$ws.Range("H31").Value = 2925.5
$ws.Range("I31").Value = 2925.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2925.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2630.5
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 2925.5
$ws.Range("I34").Value = 2925.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2925.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2723.5
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 3511.842
$ws.Range("I58").Value = 2970
$ws.Range("J58").Value = 4342.6665
$ws.Range("K58").Value = 2970
$ws.Range("L58").Value = 4342.6665
$ws.Range("M58").Value = -2767
$ws.Range("N58").Value = -4748.6665
$ws.Range("H68").Value = 41782.832
$ws.Range("J68").Value = 41782.832
$ws.Range("L68").Value = 41782.832
$ws.Range("N68").Value = -43280.832
$ws.Range("H71").Value = 41782.832
$ws.Range("J71").Value = 41782.832
$ws.Range("L71").Value = 125348.496
$ws.Range("N71").Value = -132836.496
$ws.Range("H95").Value = 23477.572
$ws.Range("J95").Value = 23477.572
$ws.Range("L95").Value = 23477.572
$ws.Range("N95").Value = -28969.572
$ws.Range("H99").Value = 30844
$ws.Range("I99").Value = 30844
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 30844
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -29346
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 9426.629000000001
$ws.Range("I122").Value = 2278.7856
$ws.Range("J122").Value = 38018
$ws.Range("K122").Value = 6836.3568
$ws.Range("L122").Value = 114054
$ws.Range("M122").Value = -4386.3568
$ws.Range("N122").Value = -118954
$ws.Range("H126").Value = 30844
$ws.Range("I126").Value = 30844
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 92532
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -90062
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 6931.684
$ws.Range("I132").Value = 7276.5884
$ws.Range("K132").Value = 21829.7652
$ws.Range("M132").Value = -19299.7652
$ws.Range("H134").Value = 3776.2
$ws.Range("I134").Value = 3673.8
$ws.Range("J134").Value = 3878.6
$ws.Range("K134").Value = 11021.4
$ws.Range("L134").Value = 11635.8
$ws.Range("M134").Value = -8486.400000000001
$ws.Range("N134").Value = -16705.8
$ws.Range("H136").Value = 3511.842
$ws.Range("I136").Value = 2970
$ws.Range("J136").Value = 4342.6665
$ws.Range("K136").Value = 8910
$ws.Range("L136").Value = 13027.9995
$ws.Range("M136").Value = -6360
$ws.Range("N136").Value = -18127.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1787850.6
$ws.Range("I4").Value = 3334031.2
$ws.Range("K4").Value = 10002093.6
$ws.Range("M4").Value = -10001981.6
$ws.Range("H9").Value = 1309875.2
$ws.Range("I9").Value = 1125146.6
$ws.Range("J9").Value = 1410636.2
$ws.Range("K9").Value = 3375439.8
$ws.Range("L9").Value = 4231908.6
$ws.Range("M9").Value = -3375215.8
$ws.Range("N9").Value = -4232356.6
$ws.Range("H39").Value = 39576.734
$ws.Range("I39").Value = 1808.5
$ws.Range("K39").Value = 5425.5
$ws.Range("M39").Value = -5131.5
$ws.Range("H60").Value = 7749.8
$ws.Range("I60").Value = 11749.667
$ws.Range("J60").Value = 1750
$ws.Range("K60").Value = 35249.001
$ws.Range("L60").Value = 5250
$ws.Range("M60").Value = -34998.001
$ws.Range("N60").Value = -5752
$ws.Range("H75").Value = 414.5
$ws.Range("J75").Value = 414.5
$ws.Range("L75").Value = 1243.5
$ws.Range("N75").Value = -3239.5
$ws.Range("H78").Value = 414.5
$ws.Range("J78").Value = 414.5
$ws.Range("L78").Value = 3730.5
$ws.Range("N78").Value = -13714.5
$ws.Range("H94").Value = 9322.875
$ws.Range("I94").Value = 3645.75
$ws.Range("K94").Value = 10937.25
$ws.Range("M94").Value = -10261.25
$ws.Range("H113").Value = 25466.791
$ws.Range("I113").Value = 408.14285
$ws.Range("K113").Value = 1224.42855
$ws.Range("M113").Value = 945.5714499999999
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H117").Value = 569.3077
$ws.Range("I117").Value = 178.33333
$ws.Range("J117").Value = 1449
$ws.Range("K117").Value = 534.99999
$ws.Range("L117").Value = 4347
$ws.Range("M117").Value = 2907.00001
$ws.Range("N117").Value = -11231
$ws.Range("H121").Value = 9972.241
$ws.Range("I121").Value = 1156.25
$ws.Range("J121").Value = 11382.8
$ws.Range("K121").Value = 3468.75
$ws.Range("L121").Value = 34148.39999999999
$ws.Range("M121").Value = -2158.75
$ws.Range("N121").Value = -36768.39999999999
$ws.Range("H123").Value = 12447.777
$ws.Range("I123").Value = 8765
$ws.Range("J123").Value = 13500
$ws.Range("K123").Value = 26295
$ws.Range("L123").Value = 40500
$ws.Range("M123").Value = -23845
$ws.Range("N123").Value = -45400
$ws.Range("H128").Value = 322768.22
$ws.Range("I128").Value = 322768.22
$ws.Range("K128").Value = 968304.6599999999
$ws.Range("M128").Value = -963324.6599999999
$ws.Range("H131").Value = 7388.778
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 29486.766
$ws.Range("I2").Value = 50039.45
$ws.Range("J2").Value = 125.78571
$ws.Range("K2").Value = 50039.45
$ws.Range("L2").Value = 125.78571
$ws.Range("M2").Value = -49926.45
$ws.Range("N2").Value = -351.78571
$ws.Range("H31").Value = 1125
$ws.Range("I31").Value = 1125
$ws.Range("K31").Value = 1125
$ws.Range("M31").Value = -833
$ws.Range("H37").Value = 1125
$ws.Range("I37").Value = 1125
$ws.Range("K37").Value = 1125
$ws.Range("M37").Value = -848
$ws.Range("H80").Value = 3613
$ws.Range("I80").Value = 3101
$ws.Range("J80").Value = 4125
$ws.Range("K80").Value = 3101
$ws.Range("L80").Value = 4125
$ws.Range("M80").Value = -2103
$ws.Range("N80").Value = -6121
$ws.Range("H83").Value = 3613
$ws.Range("I83").Value = 3101
$ws.Range("J83").Value = 4125
$ws.Range("K83").Value = 15505
$ws.Range("L83").Value = 20625
$ws.Range("M83").Value = -10513
$ws.Range("N83").Value = -30609
$ws.Range("H102").Value = 3374.4167
$ws.Range("I102").Value = 3399.3
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 3399.3
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = -1777.3
$ws.Range("N102").Value = -6494
$ws.Range("H113").Value = 1884.25
$ws.Range("I113").Value = 1914.1364
$ws.Range("K113").Value = 1914.1364
$ws.Range("M113").Value = 255.8635999999999
$ws.Range("H122").Value = 2385.3333
$ws.Range("I122").Value = 1543.1428
$ws.Range("J122").Value = 3292.3076
$ws.Range("K122").Value = 4629.428400000001
$ws.Range("L122").Value = 9876.9228
$ws.Range("M122").Value = -2179.428400000001
$ws.Range("N122").Value = -14776.9228
$ws.Range("H123").Value = 57950.668
$ws.Range("J123").Value = 57950.668
$ws.Range("L123").Value = 57950.668
$ws.Range("N123").Value = -62850.668
$ws.Range("H132").Value = 14828.08
$ws.Range("I132").Value = 15366.75
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 46100.25
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -43570.25
$ws.Range("N132").Value = -10760
$ws.Range("H136").Value = 59000
$ws.Range("J136").Value = 59000
$ws.Range("L136").Value = 177000
$ws.Range("N136").Value = -182100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1383.5518
$ws.Range("I16").Value = 1460.96
$ws.Range("K16").Value = 1460.96
$ws.Range("M16").Value = -1290.96
$ws.Range("H20").Value = 7644.4287
$ws.Range("I20").Value = 6344.6665
$ws.Range("J20").Value = 7998.909
$ws.Range("K20").Value = 6344.6665
$ws.Range("L20").Value = 7998.909
$ws.Range("M20").Value = -6118.6665
$ws.Range("N20").Value = -8450.909
$ws.Range("H39").Value = 4500
$ws.Range("I39").Value = 4500
$ws.Range("K39").Value = 4500
$ws.Range("M39").Value = -4040
$ws.Range("H40").Value = 17229.572
$ws.Range("I40").Value = 19267.834
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 19267.834
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -19131.834
$ws.Range("N40").Value = -5272
$ws.Range("H43").Value = 207375
$ws.Range("J43").Value = 207375
$ws.Range("L43").Value = 207375
$ws.Range("N43").Value = -207761
$ws.Range("H46").Value = 4076.6191
$ws.Range("I46").Value = 1154.75
$ws.Range("J46").Value = 4764.1177
$ws.Range("K46").Value = 1154.75
$ws.Range("L46").Value = 4764.1177
$ws.Range("M46").Value = -966.75
$ws.Range("N46").Value = -5140.1177
$ws.Range("H55").Value = 1371.7941
$ws.Range("I55").Value = 1210.1428
$ws.Range("K55").Value = 1210.1428
$ws.Range("M55").Value = -1037.1428
$ws.Range("H61").Value = 8198.923000000001
$ws.Range("I61").Value = 6833.5654
$ws.Range("K61").Value = 6833.5654
$ws.Range("M61").Value = -6631.5654
$ws.Range("H82").Value = 2158.8462
$ws.Range("J82").Value = 2162.5
$ws.Range("L82").Value = 2162.5
$ws.Range("N82").Value = -2884.5
$ws.Range("H85").Value = 2158.8462
$ws.Range("J85").Value = 2162.5
$ws.Range("L85").Value = 2162.5
$ws.Range("N85").Value = -4658.5
$ws.Range("H93").Value = 3595.111
$ws.Range("I93").Value = 1563.1666
$ws.Range("J93").Value = 7659
$ws.Range("K93").Value = 1563.1666
$ws.Range("L93").Value = 7659
$ws.Range("M93").Value = -315.1666
$ws.Range("N93").Value = -10155
$ws.Range("H100").Value = 3314
$ws.Range("I100").Value = 3085.3333
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3085.3333
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2544.3333
$ws.Range("N100").Value = -5082
$ws.Range("H113").Value = 8198.923000000001
$ws.Range("I113").Value = 6833.5654
$ws.Range("K113").Value = 6833.5654
$ws.Range("M113").Value = -4663.5654
$ws.Range("H122").Value = 3363.182
$ws.Range("I122").Value = 2499.1667
$ws.Range("K122").Value = 7497.500100000001
$ws.Range("M122").Value = -5047.500100000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2751.4
$ws.Range("I132").Value = 1910.6
$ws.Range("J132").Value = 5273.8
$ws.Range("K132").Value = 5731.799999999999
$ws.Range("L132").Value = 15821.4
$ws.Range("M132").Value = -3201.799999999999
$ws.Range("N132").Value = -20881.4
$ws.Range("H136").Value = 6365.731
$ws.Range("J136").Value = 11896.889
$ws.Range("L136").Value = 35690.667
$ws.Range("N136").Value = -40790.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H43").Value = 21500
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H81").Value = 7389
$ws.Range("I81").Value = 2436.25
$ws.Range("J81").Value = 23237.8
$ws.Range("K81").Value = 4872.5
$ws.Range("L81").Value = 46475.6
$ws.Range("M81").Value = -3811.5
$ws.Range("N81").Value = -48597.6
$ws.Range("H84").Value = 7389
$ws.Range("I84").Value = 2436.25
$ws.Range("J84").Value = 23237.8
$ws.Range("K84").Value = 24362.5
$ws.Range("L84").Value = 232378
$ws.Range("M84").Value = -19058.5
$ws.Range("N84").Value = -242986
$ws.Range("H122").Value = 98141.586
$ws.Range("I122").Value = 4134.75
$ws.Range("J122").Value = 286155.25
$ws.Range("K122").Value = 12404.25
$ws.Range("L122").Value = 858465.75
$ws.Range("M122").Value = -9954.25
$ws.Range("N122").Value = -863365.75
$ws.Range("H124").Value = 29960
$ws.Range("J124").Value = 29960
$ws.Range("L124").Value = 29960
$ws.Range("N124").Value = -39780
$ws.Range("H126").Value = 1859
$ws.Range("I126").Value = 1379.4
$ws.Range("J126").Value = 2098.8
$ws.Range("K126").Value = 4138.200000000001
$ws.Range("L126").Value = 6296.400000000001
$ws.Range("M126").Value = -1668.200000000001
$ws.Range("N126").Value = -11236.4
$ws.Range("H132").Value = 2625.1304
$ws.Range("I132").Value = 2032.8
$ws.Range("J132").Value = 3735.75
$ws.Range("K132").Value = 6098.4
$ws.Range("L132").Value = 11207.25
$ws.Range("M132").Value = -3568.4
$ws.Range("N132").Value = -16267.25
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H136").Value = 3740.1072
$ws.Range("I136").Value = 3740.1072
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11220.3216
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8670.321599999999
$ws.Range("N136").ClearContents()
